$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data to the table (row 48): date, minutes, description
$ws.Range("D48").Value = 44278
$ws.Range("E48").Value = 90
$ws.Range("F48").Value = "Arreglar el hilo e investigar sobre el formato de texto."

# Expand the Excel table ("Tabla1") to include the new row
$table = $ws.ListObjects.Item("Tabla1")
$table.Resize($ws.Range("D4:F48"))

# Update selection to mimic post-edit cursor position
$ws.Range("F49").Select()
